# Update "想去人数" (want-to-go count) figures in the F column across the
# relevant sheets, reflecting newly generated output data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 127
$wsExhibition.Range("F3").Value = 441

# Sheet "演出" (Performances)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 71

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 127
$wsAll.Range("F3").Value = 71
$wsAll.Range("F4").Value = 441
